$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("PLACA" stays in C, a new "destino"
# column is inserted, and the old "AJUDANTE 1"/"AJUDANTE2" columns shift
# right from D/E to E/F).
$ws.Columns("D").Insert()

# New header for the inserted column.
$ws.Range("D2").Value = "destino"

# New data point that was added alongside the new column (AJUDANTE2 for the
# first driver row).
$ws.Range("F3").Value = "MAGGIE"

# Resize the affected columns to fit their (new) contents, matching the
# widths Excel would compute automatically for the edited columns.
$ws.Columns("C:D").AutoFit()

# Restore/update the active selection to the cell the author ended up with.
$ws.Range("D3").Select() | Out-Null
